# Add more tips to the game based on player decisions.
#
# Sheet "StringLocalizations_BasicText" (4th sheet) holds key/value rows:
#   A = key, B = en-gb text, C = de, D = el, E = es (all placeholder "XXXX")
#
# Rows 81-82 are repurposed (new key + text), and 7 new rows are inserted
# below them (83-89) with the new tip strings, plus a new (empty) column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Insert 7 fresh rows at 83 (old placeholder rows get pushed down to 90/91)
# and a new column F, mirroring the structural shape of the real edit.
$ws.Rows.Item(83).Resize(7).Insert()
$ws.Columns.Item(6).Insert()

# Drop the old trailing placeholder rows that got pushed down out of the way.
$ws.Rows.Item(90).Resize(2).Delete()

# Data to write: row number, key (col A), English text (col B)
$rows = @(
    @{ R = 81; A = "TIPS_OFFICER_ONE_TURN_NEGATIVE_1";       B = "Remember thet officers sent to incidents requiring just 1 turn will return at the start of the next turn" },
    @{ R = 82; A = "TIPS_OFFICER_ONE_TURN_NEGATIVE_2";       B = "Sending officers to incidents that only require 1 turn is a great way to clear active cases" },
    @{ R = 83; A = "TIPS_OFFICER_HIGH_SEVERITY_NEGATIVE_1";  B = "Ignoring high severity cases will have a large impact on satisfaction." },
    @{ R = 84; A = "TIPS_OFFICER_HIGH_SEVERITY_NEGATIVE_2";  B = "Make sure to give high severity cases priority, failing to resolve will give a big satisfaction penalty!" },
    @{ R = 85; A = "TIPS_POSITIVE_1";                        B = "Well Done!" },
    @{ R = 86; A = "TIPS_POSITIVE_2";                        B = "Good Choice!" },
    @{ R = 87; A = "TIPS_POSITIVE_3";                        B = "Great Job!" },
    @{ R = 88; A = "TIPS_POSITIVE_4";                        B = "Awesome!" },
    @{ R = 89; A = "TIPS_POSITIVE_5";                        B = "Nice!" }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value2 = $row.A
    $ws.Cells.Item($r, 2).Value2 = $row.B
    $ws.Cells.Item($r, 3).Value2 = "XXXX"
    $ws.Cells.Item($r, 4).Value2 = "XXXX"
    $ws.Cells.Item($r, 5).Value2 = "XXXX"
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 6)).WrapText = $true
}

# New selection/scroll state matches the bottom of the (now longer) tips list.
$ws.Range("B90").Select()
